$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Starting edits"

$ws.Range("D2").Value = 2.875
$ws.Range("E2").Value = 2.921000003814697
$ws.Range("F2").Value = 3.339999914169312
$ws.Range("G2").Value = 2.684000015258789
$ws.Range("H2").Value = 445666933
$ws.Range("I2").Value = "TTD"

$ws.Range("D3").Value = 2.875
$ws.Range("E3").Value = 2.921000003814697
$ws.Range("F3").Value = 3.339999914169312
$ws.Range("G3").Value = 2.684000015258789
$ws.Range("H3").Value = 445666933
$ws.Range("I3").Value = "TTD"

$ws.Range("D4").Value = 2.875
$ws.Range("E4").Value = 2.921000003814697
$ws.Range("F4").Value = 3.339999914169312
$ws.Range("G4").Value = 2.684000015258789
$ws.Range("H4").Value = 445666933
$ws.Range("I4").Value = "TTD"

$ws.Range("D5").Value = 2.875
$ws.Range("E5").Value = 2.921000003814697
$ws.Range("F5").Value = 3.339999914169312
$ws.Range("G5").Value = 2.684000015258789
$ws.Range("H5").Value = 445666933
$ws.Range("I5").Value = "TTD"

$ws.Range("D6").Value = 2.911999940872192
$ws.Range("E6").Value = 2.517999887466431
$ws.Range("F6").Value = 2.947000026702881
$ws.Range("G6").Value = 2.349999904632568
$ws.Range("H6").Value = 445666933
$ws.Range("I6").Value = "TTD"

$ws.Range("D7").Value = 2.805000066757202
$ws.Range("E7").Value = 2.966000080108643
$ws.Range("F7").Value = 3.125
$ws.Range("G7").Value = 2.640000104904175
$ws.Range("H7").Value = 445666933
$ws.Range("I7").Value = "TTD"

$ws.Range("D8").Value = 3.719000101089478
$ws.Range("E8").Value = 3.734999895095825
$ws.Range("F8").Value = 3.844000101089478
$ws.Range("G8").Value = 3.503999948501587
$ws.Range("H8").Value = 445666933
$ws.Range("I8").Value = "TTD"

$ws.Range("D9").Value = 5.063000202178955
$ws.Range("E9").Value = 5.330999851226807
$ws.Range("F9").Value = 5.787000179290772
$ws.Range("G9").Value = 4.809000015258789
$ws.Range("H9").Value = 445666933
$ws.Range("I9").Value = "TTD"

$ws.Range("D10").Value = 6.201000213623047
$ws.Range("E10").Value = 6.592000007629395
$ws.Range("F10").Value = 6.730000019073486
$ws.Range("G10").Value = 5.85699987411499
$ws.Range("H10").Value = 445666933
$ws.Range("I10").Value = "TTD"

$ws.Range("D11").Value = 4.568999767303467
$ws.Range("E11").Value = 4.848000049591064
$ws.Range("F11").Value = 5.194900035858154
$ws.Range("G11").Value = 4.545000076293945
$ws.Range("H11").Value = 445666933
$ws.Range("I11").Value = "TTD"

$ws.Range("D12").Value = 4.895999908447266
$ws.Range("E12").Value = 5.117000102996826
$ws.Range("F12").Value = 5.500999927520752
$ws.Range("G12").Value = 4.577199935913086
$ws.Range("H12").Value = 445666933
$ws.Range("I12").Value = "TTD"

$ws.Range("D13").Value = 9.253999710083008
$ws.Range("E13").Value = 8.432000160217285
$ws.Range("F13").Value = 9.786999702453612
$ws.Range("G13").Value = 8.366000175476074
$ws.Range("H13").Value = 445666933
$ws.Range("I13").Value = "TTD"

$ws.Range("D14").Value = 15.15100002288818
$ws.Range("E14").Value = 12.35499954223633
$ws.Range("F14").Value = 15.27600002288818
$ws.Range("G14").Value = 10.5024995803833
$ws.Range("H14").Value = 445666933
$ws.Range("I14").Value = "TTD"

$ws.Range("D15").Value = 11.18599987030029
$ws.Range("E15").Value = 14.26799964904785
$ws.Range("F15").Value = 14.57800006866455
$ws.Range("G15").Value = 10.9379997253418
$ws.Range("H15").Value = 445666933
$ws.Range("I15").Value = "TTD"

$ws.Range("D16").Value = 20.1560001373291
$ws.Range("E16").Value = 22.14800071716309
$ws.Range("F16").Value = 22.63800048828125
$ws.Range("G16").Value = 19.12100028991699
$ws.Range("H16").Value = 445666933
$ws.Range("I16").Value = "TTD"

$ws.Range("D17").Value = 23.60000038146973
$ws.Range("E17").Value = 26.33099937438965
$ws.Range("F17").Value = 27.97100067138672
$ws.Range("G17").Value = 22.57999992370605
$ws.Range("H17").Value = 445666933
$ws.Range("I17").Value = "TTD"

$ws.Range("D18").Value = 18.80100059509277
$ws.Range("E18").Value = 20.07999992370605
$ws.Range("F18").Value = 21.66500091552734
$ws.Range("G18").Value = 18.03800010681152
$ws.Range("H18").Value = 445666933
$ws.Range("I18").Value = "TTD"

$ws.Range("D19").Value = 26.35199928283692
$ws.Range("E19").Value = 26.91799926757812
$ws.Range("F19").Value = 29.38999938964844
$ws.Range("G19").Value = 26.10549926757812
$ws.Range("H19").Value = 445666933
$ws.Range("I19").Value = "TTD"

$ws.Range("D20").Value = 18.48999977111816
$ws.Range("E20").Value = 29.25799942016602
$ws.Range("F20").Value = 30.60400009155273
$ws.Range("G20").Value = 15.35000038146973
$ws.Range("H20").Value = 445666933
$ws.Range("I20").Value = "TTD"

$ws.Range("D21").Value = 40.70800018310547
$ws.Range("E21").Value = 45.13199996948242
$ws.Range("F21").Value = 47.56700134277344
$ws.Range("G21").Value = 40.40000152587891
$ws.Range("H21").Value = 445666933
$ws.Range("I21").Value = "TTD"

$ws.Range("D22").Value = 52.99900054931641
$ws.Range("E22").Value = 56.64500045776367
$ws.Range("F22").Value = 67.5
$ws.Range("G22").Value = 52.60100173950195
$ws.Range("H22").Value = 445666933
$ws.Range("I22").Value = "TTD"

$ws.Range("D23").Value = 80.59999847412109
$ws.Range("E23").Value = 76.5989990234375
$ws.Range("F23").Value = 83.92800140380859
$ws.Range("G23").Value = 73.11199951171875
$ws.Range("H23").Value = 445666933
$ws.Range("I23").Value = "TTD"

$ws.Range("D24").Value = 67.39099884033203
$ws.Range("E24").Value = 72.93099975585938
$ws.Range("F24").Value = 76.86699676513672
$ws.Range("G24").Value = 63.93099975585938
$ws.Range("H24").Value = 445666933
$ws.Range("I24").Value = "TTD"

$ws.Range("D25").Value = 77.16000366210938
$ws.Range("E25").Value = 81.91000366210938
$ws.Range("F25").Value = 86.15000152587891
$ws.Range("G25").Value = 67.37000274658203
$ws.Range("H25").Value = 445666933
$ws.Range("I25").Value = "TTD"

$ws.Range("D26").Value = 70.02999877929688
$ws.Range("E26").Value = 74.91000366210938
$ws.Range("F26").Value = 82.98999786376953
$ws.Range("G26").Value = 65.31999969482422
$ws.Range("H26").Value = 445666933
$ws.Range("I26").Value = "TTD"

$ws.Range("D27").Value = 92.66000366210938
$ws.Range("E27").Value = 69.54000091552734
$ws.Range("F27").Value = 93.26000213623048
$ws.Range("G27").Value = 55.04999923706055
$ws.Range("H27").Value = 445666933
$ws.Range("I27").Value = "TTD"

$ws.Range("D28").Value = 69.84999847412109
$ws.Range("E28").Value = 58.91999816894531
$ws.Range("F28").Value = 75.98000335693359
$ws.Range("G28").Value = 56.68999862670898
$ws.Range("H28").Value = 445666933
$ws.Range("I28").Value = "TTD"

$ws.Range("D29").Value = 42.13000106811523
$ws.Range("E29").Value = 45
$ws.Range("F29").Value = 51.18999862670898
$ws.Range("G29").Value = 39
$ws.Range("H29").Value = 445666933
$ws.Range("I29").Value = "TTD"

$ws.Range("D30").Value = 60.25
$ws.Range("E30").Value = 53.2400016784668
$ws.Range("F30").Value = 64.65699768066406
$ws.Range("G30").Value = 48.15999984741211
$ws.Range("H30").Value = 445666933
$ws.Range("I30").Value = "TTD"

$ws.Range("D31").Value = 45.9900016784668
$ws.Range("E31").Value = 50.70000076293945
$ws.Range("F31").Value = 52.59999847412109
$ws.Range("G31").Value = 41.20000076293945
$ws.Range("H31").Value = 445666933
$ws.Range("I31").Value = "TTD"

$ws.Range("D32").Value = 60
$ws.Range("E32").Value = 64.33999633789062
$ws.Range("F32").Value = 65.66500091552734
$ws.Range("G32").Value = 57.16999816894531
$ws.Range("H32").Value = 445666933
$ws.Range("I32").Value = "TTD"

$ws.Range("D33").Value = 77.23999786376953
$ws.Range("E33").Value = 91.26000213623048
$ws.Range("F33").Value = 91.84999847412109
$ws.Range("G33").Value = 73.91000366210938
$ws.Range("H33").Value = 445666933
$ws.Range("I33").Value = "TTD"

$ws.Range("D34").Value = 78.15499877929688
$ws.Range("E34").Value = 70.95999908447266
$ws.Range("F34").Value = 86.42500305175781
$ws.Range("G34").Value = 64.69000244140625
$ws.Range("H34").Value = 445666933
$ws.Range("I34").Value = "TTD"

$ws.Range("D35").Value = 71.65000152587891
$ws.Range("E35").Value = 68.43000030517578
$ws.Range("F35").Value = 72.15000152587891
$ws.Range("G35").Value = 61.47499847412109
$ws.Range("H35").Value = 445666933
$ws.Range("I35").Value = "TTD"

$ws.Range("D36").Value = 87.40000152587891
$ws.Range("E36").Value = 82.84999847412109
$ws.Range("F36").Value = 88.61000061035156
$ws.Range("G36").Value = 76.12000274658203
$ws.Range("H36").Value = 445666933
$ws.Range("I36").Value = "TTD"

$ws.Range("D37").Value = 97.43000030517578
$ws.Range("E37").Value = 89.87999725341797
$ws.Range("F37").Value = 102.6699981689453
$ws.Range("G37").Value = 88.26000213623047
$ws.Range("H37").Value = 445666933
$ws.Range("I37").Value = "TTD"

$ws.Range("D38").Value = 110.5849990844727
$ws.Range("E38").Value = 120.2099990844727
$ws.Range("F38").Value = 123.8499984741211
$ws.Range("G38").Value = 106.6999969482422
$ws.Range("H38").Value = 445666933
$ws.Range("I38").Value = "TTD"

$ws.Range("D39").Value = 119.0699996948242
$ws.Range("E39").Value = 118.6800003051758
$ws.Range("F39").Value = 127.5899963378906
$ws.Range("G39").Value = 115.8499984741211
$ws.Range("H39").Value = 445666933
$ws.Range("I39").Value = "TTD"

$ws.Range("D40").Value = 54.54000091552734
$ws.Range("E40").Value = 53.63000106811523
$ws.Range("F40").Value = 57.67800140380859
$ws.Range("G40").Value = 42.95999908447266
$ws.Range("H40").Value = 445666933
$ws.Range("I40").Value = "TTD"

$ws.Range("D41").Value = 73.43499755859375
$ws.Range("E41").Value = 86.95999908447266
$ws.Range("F41").Value = 89.13400268554688
$ws.Range("G41").Value = 72.37400054931641
$ws.Range("H41").Value = 445666933
$ws.Range("I41").Value = "TTD"
Write-Host "Edits complete"
